$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Quandl Merrill Lynch bond index meta-feature rows to append below
# the existing data (rows 2-31 already populated; next free row is 32).
$newRows = @(
    @{ Key = "ML/EMHYY";  Expl = "EM_HY_YIELD" },
    @{ Key = "ML/AAAEY";  Expl = "US_IG_YIELD" },
    @{ Key = "ML/BBBEY";  Expl = "US_CORP_YIELD" },
    @{ Key = "ML/EMCTRI"; Expl = "EM_CORP_RET" },
    @{ Key = "ML/USTRI";  Expl = "US_HY_YIELD" },
    @{ Key = "ML/EMHGY";  Expl = "EM_IG_YIELD" },
    @{ Key = "ML/EMHG";   Expl = "EMEA_CORP_YIELD" }
)

$startRow = 32
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = "QUANDL"      # SOURCE
    $ws.Cells.Item($r, 2).Value = $row.Key       # KEY
    $ws.Cells.Item($r, 3).Value = $row.Expl      # EXPL
    $ws.Cells.Item($r, 4).Value = "LAST"         # NA_METHOD

    # Writing the literal text "FALSE" via .Value auto-coerces to a real
    # Boolean (matches Excel's native typed-entry behaviour). The source
    # column stores it as literal text (matching the other HLOC cells in
    # the sheet), so round-trip it through a formula -> values-only paste
    # to land a plain shared-string "FALSE" instead of a boolean cell.
    $gCell = $ws.Cells.Item($r, 7)               # HLOC
    $gCell.Formula = '="FALSE"'
    $gCell.Copy()
    $gCell.PasteSpecial(-4163)                    # xlPasteValues

    $ws.Cells.Item($r, 8).Value = "B"            # FREQ
}
$excel.CutCopyMode = 0

# Reflect the new used range / viewport like the saved workbook did:
# scrolled so row 17 is at the top, with the newly added rows selected.
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 17
$null = $ws.Range("A31:A38").Select()
